$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.624.09'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').Value = '''1.643.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  +0.69%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.26%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''215.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '''  +1.04%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.01%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.85%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''19.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  +0.30%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.02%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''1.872.35'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  +0.66%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = '''Polkadot'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = '''4.21'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  +3.19%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = '''WrappedEther'
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = '''1.623.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.27%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.50%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''65.88'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +4.26%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''26.671.60'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  +0.09%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.0₃0750'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  +1.44%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''218.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.09%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''1.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.27%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = '''  +2.41%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  +2.29%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''9.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  +1.93%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''2.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  +9.27%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''146.28'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -1.14%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  +0.33%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E28').Value = '''  +3.90%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = '''  +2.33%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''0.0518'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +2.81%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.58%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  +3.16%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = '''  +2.65%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''1.275.10'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  +5.52%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  +2.71%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '''  +5.97%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  +0.06%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = '''  +5.97%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.829'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +2.67%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''1.01'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  +0.25%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = '''  +2.15%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = '''  -1.46%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = '''  +1.38%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''1.783.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +0.68%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''93.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  +0.13%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = '''59.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +9.33%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = '''  +3.72%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = '''  +0.84%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''7.78'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  +2.21%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  +4.19%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = '''  -0.44%  '
$ws.Range('E51').Style = 'Normal'

Write-Host "Applied 73 cell updates"